# The commit swaps the two embedded theme parts: the slide master's theme
# (ppt/theme/theme1.xml, previously the "Integral" / "Red Violet" theme)
# ends up holding the stock "Office Theme" / "Office" color scheme, while
# the notes-master theme (ppt/theme/theme2.xml) ends up holding the
# "Integral" / "Red Violet" colors that used to live in theme1.xml.
#
# Font scheme (majorFont/minorFont) and format scheme (fills/lines/effects)
# are already byte-identical between the two theme parts, so the only
# observable change is the 12-slot color scheme (clrScheme) used by the
# slide master / slides, which is reachable and settable through
# Master.Theme.ThemeColorScheme.Colors(n).RGB.
#
# COM RGB values are packed as 0x00BBGGRR, so each target hex color
# AABBCC is written as (CC*65536)+(BB*256)+AA.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# 1 dk1      -> 000000
$colorScheme.Colors(1).RGB = 0
# 2 lt1      -> FFFFFF
$colorScheme.Colors(2).RGB = 16777215
# 3 dk2      -> 44546A
$colorScheme.Colors(3).RGB = 6968388
# 4 lt2      -> E7E6E6
$colorScheme.Colors(4).RGB = 15132391
# 5 accent1  -> 5B9BD5
$colorScheme.Colors(5).RGB = 13998939
# 6 accent2  -> ED7D31
$colorScheme.Colors(6).RGB = 3243501
# 7 accent3  -> A5A5A5
$colorScheme.Colors(7).RGB = 10855845
# 8 accent4  -> FFC000
$colorScheme.Colors(8).RGB = 49407
# 9 accent5  -> 4472C4
$colorScheme.Colors(9).RGB = 12874308
# 10 accent6 -> 70AD47
$colorScheme.Colors(10).RGB = 4697456
# 11 hlink   -> 0563C1
$colorScheme.Colors(11).RGB = 12673797
# 12 folHlink -> 954F72
$colorScheme.Colors(12).RGB = 7491477
